$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix company/name strings: stray commas used where a period was intended ---
$nameFixes = @{
    "E21"  = "FERNANDEZ MARIO H. GALLICET OSCAR M"
    "E68"  = "FERNANDEZ MARIO H. GALLICET OSCAR M"
    "E69"  = "IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA"
    "F69"  = "IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA"
    "E70"  = "MARSICO GUILLERMO MIGUEL. MARSICO JUAN EDUARDO"
    "E92"  = "TRABICHET MARIA. VERGARA ADEL Y OTRA"
    "F92"  = "TRABICHET MARIA. VERGARA ADEL Y OTRA"
    "E136" = "MARSICO GUILLERMO MIGUEL. MARSICO JUAN EDUARDO"
}

foreach ($addr in $nameFixes.Keys) {
    $ws.Range($addr).Value = $nameFixes[$addr]
}

# --- Fix "Importe" (column H) values: scraped with Spanish-locale formatting
#     (", " decimal separator, "." thousands separator) -- convert to plain
#     decimal-point strings, e.g. "2.787,50" -> "2787.50", keeping them as
#     TEXT (not real numbers) exactly as in the source data. ---
$importeValues = @(
    "2787.50",
    "1349.35",
    "5209.60",
    "1907.50",
    "131368.16",
    "22580.50",
    "5856.00",
    "15452.75",
    "2045.44",
    "8684.97",
    "1534.50",
    "8287.85",
    "2400.00",
    "2400.00",
    "1780.00",
    "9.59",
    "48.86",
    "45.50",
    "28977.70",
    "280.00",
    "13379.65",
    "99.42",
    "115.00",
    "1317.47",
    "40.00",
    "232.00",
    "13104.90",
    "186.80",
    "1724.24",
    "945.00",
    "1705.13",
    "95.00",
    "650.00",
    "1464.00",
    "2140.00",
    "1430.00",
    "94.00",
    "2.66",
    "2287.82",
    "1092.17",
    "8786.65",
    "467.50",
    "1064.00",
    "6750.00",
    "307.59",
    "120.00",
    "1100.00",
    "386.90",
    "1003.24",
    "69.00",
    "1200.00",
    "15.50",
    "271.76",
    "124.00",
    "2088.70",
    "630.00",
    "1380.00",
    "612.20",
    "200.00",
    "22560.00",
    "2302.00",
    "1200.00",
    "754.50",
    "2970.00",
    "268.00",
    "1337.00",
    "753.00",
    "243.43",
    "1415.00",
    "7835.00",
    "162.00",
    "4663.90",
    "1450.00",
    "3000.00",
    "144870.00",
    "7040.04",
    "190.78",
    "120.00",
    "7.88",
    "30.00",
    "46.00",
    "3393.10",
    "2395.00",
    "2094.50",
    "110.00",
    "311.00",
    "42.50",
    "250.00",
    "115.00",
    "164.00",
    "49.00",
    "1186.10",
    "1795.10",
    "613.92",
    "219.00",
    "223.99",
    "5591.00",
    "1650.00",
    "10090.00",
    "37.00",
    "2337.00",
    "349.00",
    "3322.00",
    "15785.01",
    "2210.00",
    "3025.00",
    "2245.00",
    "375.00",
    "16868.00",
    "4900.00",
    "2921.00",
    "785.00",
    "500.00",
    "3500.00",
    "1690.86",
    "11572.50",
    "1700.90",
    "420.60",
    "470.95",
    "2900.00",
    "6990.00",
    "1400.00",
    "1000.00",
    "773.50",
    "250.00",
    "700.00",
    "16717.90",
    "200.00",
    "4250.00",
    "1500.00",
    "20660.00",
    "3000.00",
    "249.90",
    "34545.50",
    "20.00",
    "2461.00",
    "1600.00",
    "9980.00",
    "736.92",
    "1955.00",
    "4856.00",
    "360.00",
    "131.40",
    "111.41",
    "960.00",
    "250.34",
    "6780.63",
    "3563.59",
    "5372.00",
    "2312.70",
    "3843.20",
    "1950.00",
    "8475.69",
    "3975.00",
    "98.76",
    "885.60",
    "582.40",
    "2718.16",
    "3872.61",
    "2791.00",
    "13100.00",
    "224900.00",
    "189854.50",
    "243375.00",
    "222700.00",
    "249084.00",
    "245924.00",
    "200900.00",
    "186556.46",
    "2400.00"
)

$rng = $ws.Range("H2:H171")
$rng.NumberFormat = "@"
$data = New-Object 'object[,]' $importeValues.Length,1
for ($i = 0; $i -lt $importeValues.Length; $i++) {
    $data[$i,0] = $importeValues[$i]
}
$rng.Value = $data
$rng.Style = "Normal"
